$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row for the Skellam "Mean-Dispersion" parametrization ---
# The new row goes right after the existing "Skellam / Difference" row (row 20),
# i.e. becomes the new row 21, pushing every row below it down by one.
$ws.Rows.Item(21).Insert()

# Seed the new row from the row above (Skellam/Difference) so that text vs.
# boolean-like ("TRUE"/"FALSE") cell types and styles are copied verbatim,
# instead of being re-inferred from literal strings.
$ws.Range("A20:H20").Copy($ws.Range("A21:H21"))

# Now overwrite only the two cells that actually differ for this parametrization.
$ws.Range("B21").Value = "Mean-Dispersion"
$ws.Range("D21").Value = "meandisp"

# --- Refresh the AutoFilter so its range covers the grown table ---
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}
$ws.Range("A1:H28").AutoFilter()

# --- Update the (LibreOffice-style) hidden filter-database defined names ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "dist_table!_FilterDatabase") {
        $n.RefersTo = "=dist_table!`$A`$1:`$H`$28"
    }
    if ($n.Name -eq "dist_table!_FilterDatabase_0_0") {
        $n.RefersTo = "=dist_table!`$A`$1:`$H`$26"
    }
}

# --- Restore the active selection to H21, matching the saved view state ---
$ws.Range("H21").Select()
